$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Ranking")
$ws.Range("H6").Value = 0.01256259791707335
$ws.Range("I6").Value = 0.009679735796162117
$ws.Range("H7").Value = 0.01173646236175684
$ws.Range("I7").Value = 0.009338790450295106
$ws.Range("H9").Value = 0.06498615654567919
$ws.Range("I9").Value = 0.06600441838929333

$ws = $wb.Worksheets.Item("Matriz_Pvalores")
$ws.Range("G2").Value = 0.0000478992950014856
$ws.Range("H2").Value = 0.00004877003845815686
$ws.Range("J2").Value = 0.001427792481607959
$ws.Range("G3").Value = 0.006569872598065851
$ws.Range("H3").Value = 0.009736967387809159
$ws.Range("J3").Value = 0.0000004014213428327196
$ws.Range("G4").Value = 0.01139646693791385
$ws.Range("H4").Value = 0.02122358983114836
$ws.Range("J4").Value = 0.0000004306503971207576
$ws.Range("G5").Value = 0.9661412312661364
$ws.Range("H5").Value = 0.8608817972672327
$ws.Range("J5").Value = 0.00000001163064689002624
$ws.Range("G6").Value = 0.00007719447554777048
$ws.Range("H6").Value = 0.00007871524290248288
$ws.Range("J6").Value = 0.4922838102681615
$ws.Range("B7").Value = 0.0000478992950014856
$ws.Range("C7").Value = 0.006569872598065851
$ws.Range("D7").Value = 0.01139646693791385
$ws.Range("E7").Value = 0.9661412312661364
$ws.Range("F7").Value = 0.00007719447554777048
$ws.Range("H7").Value = 0.6003746506766476
$ws.Range("I7").Value = 0.0118974966548806
$ws.Range("J7").Value = 0.0000335641213466964
$ws.Range("B8").Value = 0.00004877003845815686
$ws.Range("C8").Value = 0.009736967387809159
$ws.Range("D8").Value = 0.02122358983114836
$ws.Range("E8").Value = 0.8608817972672327
$ws.Range("F8").Value = 0.00007871524290248288
$ws.Range("G8").Value = 0.6003746506766476
$ws.Range("I8").Value = 0.01685804575200978
$ws.Range("J8").Value = 0.00003483864252684654
$ws.Range("G9").Value = 0.0118974966548806
$ws.Range("H9").Value = 0.01685804575200978
$ws.Range("J9").Value = 0.0000002978936606279348
$ws.Range("B10").Value = 0.001427792481607959
$ws.Range("C10").Value = 0.0000004014213428327196
$ws.Range("D10").Value = 0.0000004306503971207576
$ws.Range("E10").Value = 0.00000001163064689002624
$ws.Range("F10").Value = 0.4922838102681615
$ws.Range("G10").Value = 0.0000335641213466964
$ws.Range("H10").Value = 0.00003483864252684654
$ws.Range("I10").Value = 0.0000002978936606279348

$ws = $wb.Worksheets.Item("Matriz_DM_Original")
$ws.Range("G2").Value = 10.33817915900059
$ws.Range("H2").Value = 10.30566826406916
$ws.Range("J2").Value = 5.563768475740805
$ws.Range("G3").Value = -4.07080263500573
$ws.Range("H3").Value = -3.7299740552165
$ws.Range("J3").Value = -23.38171515588723
$ws.Range("G4").Value = -3.597807817614061
$ws.Range("H4").Value = -3.09603618672405
$ws.Range("J4").Value = -23.10678826514726
$ws.Range("G5").Value = 0.04424975004056525
$ws.Range("H5").Value = 0.1829245400964258
$ws.Range("J5").Value = -42.32998546331976
$ws.Range("G6").Value = 9.50747858442066
$ws.Range("H6").Value = 9.474831469846844
$ws.Range("J6").Value = -0.7310756733642315
$ws.Range("B7").Value = -10.33817915900059
$ws.Range("C7").Value = 4.07080263500573
$ws.Range("D7").Value = 3.597807817614061
$ws.Range("E7").Value = -0.04424975004056525
$ws.Range("F7").Value = -9.50747858442066
$ws.Range("H7").Value = 0.5527984338458297
$ws.Range("I7").Value = 3.562063661480992
$ws.Range("J7").Value = -10.99940859250579
$ws.Range("B8").Value = -10.30566826406916
$ws.Range("C8").Value = 3.7299740552165
$ws.Range("D8").Value = 3.09603618672405
$ws.Range("E8").Value = -0.1829245400964258
$ws.Range("F8").Value = -9.474831469846844
$ws.Range("G8").Value = -0.5527984338458297
$ws.Range("I8").Value = 3.278314725130754
$ws.Range("J8").Value = -10.92834601486896
$ws.Range("G9").Value = -3.562063661480992
$ws.Range("H9").Value = -3.278314725130754
$ws.Range("J9").Value = -24.58460814617296
$ws.Range("B10").Value = -5.563768475740805
$ws.Range("C10").Value = 23.38171515588723
$ws.Range("D10").Value = 23.10678826514726
$ws.Range("E10").Value = 42.32998546331976
$ws.Range("F10").Value = 0.7310756733642315
$ws.Range("G10").Value = 10.99940859250579
$ws.Range("H10").Value = 10.92834601486896
$ws.Range("I10").Value = 24.58460814617296
